# Update countries & provincias Spain
# - Reorder two pairs of countries in the "Pais" sheet (their case-count
#   rows keep their position/rank, but the country label assigned to a
#   couple of rows swaps, matching the shared-string reshuffle in the
#   diff):
#     * row 143..146: Uruguay/Jordania/Malta/Bahamas -> Bahamas/Uruguay/Jordania/Malta
#     * row 213..214: Islas Malvinas/Montserrat -> Montserrat/Islas Malvinas
# - Refresh the day's case counters for a handful of rows.
# - Refresh the "Datos actualizados a ..." timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country label swaps -------------------------------------------------
$ws.Range("A143").Value = "Bahamas"
$ws.Range("A144").Value = "Uruguay"
$ws.Range("A145").Value = "Jordania"
$ws.Range("A146").Value = "Malta"

$ws.Range("A213").Value = "Montserrat"
$ws.Range("A214").Value = "Islas Malvinas"

# --- Updated per-country counters ----------------------------------------
# row => @{ column letter = new value }
$rowUpdates = @{
    4   = @{ B = 5699221; C = 43247; D = 3060534; E = 2462385;             G = 1228; H = 176302 }
    27  = @{ B = 123490;  C = 336;   D = 109822;  E = 4619 }
    114 = @{ B = 4174;    C = 42;    D = 3127;    E = 967 }
    143 = @{ B = 1531;    C = 107;   D = 209;     E = 1300;              G = 2;    H = 22 }
    144 = @{ B = 1493;    C = 8;     D = 1228;    E = 225;                         H = 40 }
    145 = @{ B = 1482;    C = 44;    D = 1259;    E = 212;                         H = 11 }
    146 = @{ B = 1470;    C = 47;    D = 784;     E = 677;                         H = 9 }
    186 = @{ B = 205;     C = 2;                  E = 2 }
    196 = @{                        D = 56;      E = 2 }
    213 = @{                        D = 12;                                       H = 1 }
    214 = @{                        D = 13;                                       H = 0 }
}

foreach ($r in $rowUpdates.Keys) {
    $cols = $rowUpdates[$r]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$r").Value = $cols[$col]
    }
}

# --- Timestamp -------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Agosto de 2020 a las 02:57"
